# FighterSkills.xlsx edit
#
# - Adds the missing "[[MP: 1 ]]" cost line to the "10-LiquidKick" skill
#   sheet (a new row inserted right under the "[[AP: 4 ]]" line, mirroring
#   the layout already used on the "7-EagleFlight" sheet), which pushes
#   every row below it down by one.
# - Minor navigation/selection state left behind by the edit session:
#   the "8-IronPuch" sheet ends up with its selection on E10, the
#   "9-KiBlast" sheet is scrolled down a bit, and the final active sheet
#   (and tab) is "10-LiquidKick" with the newly edited cell selected.

$wb = $excel.ActiveWorkbook

# --- 8-IronPuch: selection left on E10 -------------------------------
$wsIronPunch = $wb.Worksheets.Item("8-IronPuch")
$wsIronPunch.Activate() | Out-Null
$wsIronPunch.Range("E10").Select() | Out-Null

# --- 9-KiBlast: view scrolled down a little ---------------------------
$wsKiBlast = $wb.Worksheets.Item("9-KiBlast")
$wsKiBlast.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 4

# --- 10-LiquidKick: add the new [[MP: 1 ]] line -----------------------
$wsLiquidKick = $wb.Worksheets.Item("10-LiquidKick")
$wsLiquidKick.Activate() | Out-Null

# Insert a new row above the old row 8 ("[[Range: 1 ]]"), shifting the
# rest of the skill card down by one row.
$wsLiquidKick.Rows.Item(8).Insert() | Out-Null
$wsLiquidKick.Range("C8").Value = "[[MP: 1 ]]"
$wsLiquidKick.Rows.Item(8).RowHeight = 19.5

$wsLiquidKick.Range("C29").Select() | Out-Null
